# Before Class Week 2
# Insert a new worksheet "3.1" before the existing "4.4" sheet and populate
# it with precipitation data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "3.1"

$ws.Range("A1").Value = "Precip"
$ws.Range("A2").Value = 26.87
$ws.Range("A3").Value = 26.94
$ws.Range("A4").Value = 28.28
$ws.Range("A5").Value = 29.48
$ws.Range("A6").Value = 31.56
$ws.Range("A7").Value = 32.78
$ws.Range("A8").Value = 33.07
$ws.Range("A9").Value = 33.62
$ws.Range("A10").Value = 34.98
$ws.Range("A11").Value = 35.09
$ws.Range("A12").Value = 35.2
$ws.Range("A13").Value = 35.38
$ws.Range("A14").Value = 35.96
$ws.Range("A15").Value = 36.02
$ws.Range("A16").Value = 36.65
$ws.Range("A17").Value = 36.83
$ws.Range("A18").Value = 36.99
$ws.Range("A19").Value = 38.15
$ws.Range("A20").Value = 39.34
$ws.Range("A21").Value = 39.62
$ws.Range("A22").Value = 39.86
$ws.Range("A23").Value = 40.21
$ws.Range("A24").Value = 40.54
$ws.Range("A25").Value = 41.11
$ws.Range("A26").Value = 41.34
$ws.Range("A27").Value = 41.44
$ws.Range("A28").Value = 41.46
$ws.Range("A29").Value = 41.94
$ws.Range("A30").Value = 43.3
$ws.Range("A31").Value = 43.53
$ws.Range("A32").Value = 45.62
$ws.Range("A33").Value = 46.02
$ws.Range("A34").Value = 47.73
$ws.Range("A35").Value = 47.9
$ws.Range("A36").Value = 48.02
$ws.Range("A37").Value = 50.5
$ws.Range("A38").Value = 51.17
$ws.Range("A39").Value = 51.97
$ws.Range("A40").Value = 54.29
$ws.Range("A41").Value = 57.54

$ws.Range("A42").Select() | Out-Null
